$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(68).Copy()
$ws.Rows.Item(69).Insert(-4121)
